$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A values (rows 2-11) from 1 to 4
$ws.Range("A2:A11").Value = 4

# Update the selected cell/range to J7
$ws.Range("J7").Select()
